$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "language" column K ---
$ws.Columns("K").Insert()
$ws.Range("K1").Value = "language"

# --- Fix row 5: D5 and H5 should become real numbers (previously stored as text) ---
$ws.Range("D5").Value = 508988751
$ws.Range("H5").Value = 5

# --- Append new rows 6-9 ---
$ws.Range("A6").Value = "ORD-20250301071419"
$ws.Range("B6").Value = "2025-03-01 07:14:19"
$ws.Range("C6").Value = "Steven"
$ws.Range("D6").Value = 508988751
$ws.Range("E6").Value = "Jeans"
$ws.Range("F6").Value = "S"
$ws.Range("G6").Value = "Green"
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = "Abu Dhabi"
$ws.Range("J6").Value = "New"

$ws.Range("A7").Value = "ORD-20250301073023"
$ws.Range("B7").Value = "2025-03-01 07:30:23"
$ws.Range("C7").Value = "Can i see the product"
$ws.Range("D7").Value = "My name is fee sidabalok"
$ws.Range("E7").Value = "I want to order dress"
$ws.Range("F7").Value = "M"
$ws.Range("G7").Value = "Red"
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = "Jalan aries utama no.57 rt.3/rw.3 Meruya utara kembangan jakarta barat"
$ws.Range("J7").Value = "New"

$ws.Range("A8").Value = "ORD-20250301172534"
$ws.Range("B8").Value = "2025-03-01 17:25:34"
$ws.Range("C8").Value = "Steven"
$ws.Range("D8").Formula = '="0508988751"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "Rok"
$ws.Range("F8").Value = "S"
$ws.Range("G8").Value = "Ungu"
$ws.Range("H8").Value = 5
$ws.Range("I8").Value = "Abu dhabi"
$ws.Range("J8").Value = "New"
$ws.Range("K8").Value = "id"

$ws.Range("A9").Value = "ORD-20250301172706"
$ws.Range("B9").Value = "2025-03-01 17:27:06"
$ws.Range("C9").Value = "Steven"
$ws.Range("D9").Formula = '="0508988751"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "Jeans"
$ws.Range("F9").Value = "XXL"
$ws.Range("G9").Value = "White"
$ws.Range("H9").Formula = '="5"'
$ws.Range("H9").Copy()
$ws.Range("H9").PasteSpecial(-4163)
$ws.Range("I9").Value = "Abu Dhabi"
$ws.Range("J9").Value = "New"
$ws.Range("K9").Value = "en"
